$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (16 and 17) had their "Periodo Mora" (column E) and
# "Valor Mora" (column F) values swapped between them:
#   Row16: Periodo "2402" -> "2401", Valor Mora 64000 -> 40533
#   Row17: Periodo "2401" -> "2402", Valor Mora 40533 -> 64000
$ws.Range("E16").Value = "2401"
$ws.Range("F16").Value = 40533
$ws.Range("E17").Value = "2402"
$ws.Range("F17").Value = 64000
